$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each division problem lives in its own table cell. Several problems share the
# same text at different points in the edit (e.g. "49÷6=" appears twice, and a
# couple of the new values collide with other cells' old values mid-sequence),
# so every replacement below is scoped to a single cell's Range and uses
# wdReplaceOne (the final "1" argument) instead of wdReplaceAll to guarantee it
# only ever touches that one cell's run.

$cell = $t.Cell(1, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("22÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷3=", 1) | Out-Null

$cell = $t.Cell(1, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("68÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷2=", 1) | Out-Null

$cell = $t.Cell(1, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("53÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷7=", 1) | Out-Null

$cell = $t.Cell(1, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("43÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷7=", 1) | Out-Null

$cell = $t.Cell(1, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("16÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷9=", 1) | Out-Null

$cell = $t.Cell(5, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("40÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷5=", 1) | Out-Null

$cell = $t.Cell(5, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("22÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷2=", 1) | Out-Null

$cell = $t.Cell(5, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("49÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷5=", 1) | Out-Null

$cell = $t.Cell(5, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("52÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷7=", 1) | Out-Null

$cell = $t.Cell(5, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("31÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷8=", 1) | Out-Null

$cell = $t.Cell(9, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("39÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷8=", 1) | Out-Null

$cell = $t.Cell(9, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("73÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷2=", 1) | Out-Null

$cell = $t.Cell(9, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("69÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷9=", 1) | Out-Null

$cell = $t.Cell(9, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("97÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷8=", 1) | Out-Null

$cell = $t.Cell(9, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("10÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷2=", 1) | Out-Null

$cell = $t.Cell(13, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("36÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷5=", 1) | Out-Null

$cell = $t.Cell(13, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("92÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷9=", 1) | Out-Null

$cell = $t.Cell(13, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("18÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷7=", 1) | Out-Null

$cell = $t.Cell(13, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("60÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷8=", 1) | Out-Null

$cell = $t.Cell(13, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("49÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷5=", 1) | Out-Null

$cell = $t.Cell(17, 1)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("14÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷2=", 1) | Out-Null

$cell = $t.Cell(17, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("94÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷8=", 1) | Out-Null

$cell = $t.Cell(17, 3)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("49÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷7=", 1) | Out-Null

$cell = $t.Cell(17, 4)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("47÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "91÷6=", 1) | Out-Null

$cell = $t.Cell(17, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Find.Execute("93÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷7=", 1) | Out-Null
